$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Coin name and Link for rows 10 and 11 (OKB <-> Polygon)
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"

# Update Price (D) and Volume(1h) (E) columns for each row
# Price column is stored as text; force text format for numeric-looking
# values so Excel does not auto-convert them to numbers and lose exact
# formatting (e.g. trailing zeros).
$ws.Range("D2").Value = "23.865.29"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").Value = "1.620.74"
$ws.Range("E3").Value = "  -3.10%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.97"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3936"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3845"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9992"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.365"
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "49.67"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08461"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.91"
$ws.Range("E13").Value = "  -5.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.060"
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.576"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001280"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "1.613.15"
$ws.Range("E17").Value = "  -4.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.00"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06926"
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.10"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.816"
$ws.Range("E21").Value = "  -3.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.42"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").Value = "23.864.06"
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.461"
$ws.Range("E25").Value = "  +4.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.863"
$ws.Range("E26").Value = "  +2.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.26"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.54"
$ws.Range("E28").Value = "  -2.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.93"
$ws.Range("E29").Value = "  -3.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.274"
$ws.Range("E30").Value = "  -9.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.869"
$ws.Range("E31").Value = "  -5.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.505"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "1.791.55"
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08132"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9796"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02891"
$ws.Range("E36").Value = "  -6.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.604"
$ws.Range("E37").Value = "  -5.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2672"
$ws.Range("E38").Value = "  -5.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09153"
$ws.Range("E39").Value = "  -4.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.40"
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.61"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.422"
$ws.Range("E42").Value = "  -6.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7517"
$ws.Range("E43").Value = "  -4.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.10"
$ws.Range("E44").Value = "  -3.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6913"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.476"
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.071"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9995"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08247"
$ws.Range("E49").Value = "  -4.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.06"
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.210"
$ws.Range("E51").Value = "  -8.86%  "
